# Insert a new data row at row 255 of the (single-sheet) price-list table.
# All rows currently at 255..312 shift down to 256..313, and the sheet's
# used-range dimension grows from A1:R312 to A1:R313 automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(255).Insert()

$ws.Cells.Item(255, 1).Value  = 10
$ws.Cells.Item(255, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(255, 3).Value  = "La Araucanía"
$ws.Cells.Item(255, 4).Value  = 45211
$ws.Cells.Item(255, 5).Value  = 9
$ws.Cells.Item(255, 6).Value  = 100112012
$ws.Cells.Item(255, 7).Value  = "Espinaca"
$ws.Cells.Item(255, 8).Value  = "Sin especificar"
$ws.Cells.Item(255, 9).Value  = "Primera"
$ws.Cells.Item(255, 10).Value = 90
$ws.Cells.Item(255, 11).Value = 10000
$ws.Cells.Item(255, 12).Value = 10000
$ws.Cells.Item(255, 13).Value = 10000
$ws.Cells.Item(255, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(255, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(255, 16).Value = 833
$ws.Cells.Item(255, 17).Value = 12
$ws.Cells.Item(255, 18).Value = "Hortaliza"
